$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.784.34'
$ws.Range('E2').Value = '  +0.79%  '
$ws.Range('D3').Value = '3.582.21'
$ws.Range('E3').Value = '  +1.14%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '603.75'
$ws.Range('E5').Value = '  +1.14%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '137.36'
$ws.Range('E6').Value = '  -0.92%  '
$ws.Range('D7').Value = '3.581.42'
$ws.Range('E7').Value = '  +1.05%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  +1.09%  '
$ws.Range('E10').Value = '  +0.41%  '
$ws.Range('E11').Value = '  +5.37%  '
$ws.Range('E12').Value = '  +1.23%  '
$ws.Range('D13').Value = '4.190.42'
$ws.Range('E13').Value = '  +1.15%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '28.29'
$ws.Range('E14').Value = '  +3.61%  '
$ws.Range('E15').Value = '  +0.00%  '
$ws.Range('D16').Value = '3.579.74'
$ws.Range('E16').Value = '  +1.00%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.116'
$ws.Range('E17').Value = '  -0.29%  '
$ws.Range('D18').Value = '65.869.61'
$ws.Range('E18').Value = '  +0.90%  '
$ws.Range('E19').Value = '  -2.18%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.68'
$ws.Range('E20').Value = '  +2.04%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.87'
$ws.Range('E21').Value = '  -1.56%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '394.84'
$ws.Range('E22').Value = '  +0.18%  '
$ws.Range('E23').Value = '  +2.48%  '
$ws.Range('D24').Value = '3.728.56'
$ws.Range('E24').Value = '  +1.28%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '74.15'
$ws.Range('E25').Value = '  +0.33%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').Value = '  +1.62%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.12'
$ws.Range('E28').Value = '  +3.67%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.62'
$ws.Range('E29').Value = '  +25.99%  '
$ws.Range('E30').Value = '  +2.58%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '8.59'
$ws.Range('E31').Value = '  +5.24%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.996'
$ws.Range('E32').Value = '  -0.24%  '
$ws.Range('D33').Value = '3.589.33'
$ws.Range('E33').Value = '  +1.07%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '24.53'
$ws.Range('E34').Value = '  +2.88%  '
$ws.Range('E35').Value = '  +1.92%  '
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.38'
$ws.Range('E37').Value = '  +7.46%  '
$ws.Range('E38').Value = '  +4.80%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '7.02'
$ws.Range('E39').Value = '  +0.48%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '168.23'
$ws.Range('E40').Value = '  -0.66%  '
$ws.Range('E41').Value = '  +4.07%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.838'
$ws.Range('E42').Value = '  +1.61%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '26.76'
$ws.Range('E43').Value = '  +1.54%  '
$ws.Range('E44').Value = '  +6.67%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '43.08'
$ws.Range('E45').Value = '  +0.68%  '
$ws.Range('E46').Value = '  +2.37%  '
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.70'
$ws.Range('E48').Value = '  +1.37%  '
$ws.Range('E49').Value = '  +2.93%  '
$ws.Range('D50').Value = '2.462.96'
$ws.Range('E50').Value = '  +3.29%  '
$ws.Range('E51').Value = '  +4.55%  '
